
$wb = $excel.ActiveWorkbook

# F-column ("want to go" count) updates keyed by the FINAL row number (after the two row
# deletions below). The same row numbers / values apply identically on both the "展览" and
# "全部类型" sheets, since they enumerate the same events.
$fUpdates = @{
    2  = 282
    3  = 13719
    4  = 1343
    7  = 178
    8  = 261
    9  = 488
    10 = 11
    11 = 82
    15 = 442
    16 = 5683
    17 = 121
    18 = 78
    20 = 61
    22 = 141
    23 = 196
}

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "展览" -or $ws.Name -eq "全部类型") {
        # Two events have dropped out of the source feed:
        #   row 7 = "巢湖·原铁崩ONLY"
        #   row 3 = "合肥·第十三届合肥次元之门动漫游戏博览会-多多poi&Mace专场"
        # Delete the higher-numbered row first so the second delete's row index stays valid.
        $ws.Rows.Item(7).Delete()
        $ws.Rows.Item(3).Delete()

        $lastRow = $ws.UsedRange.Rows.Count

        # Column A holds a plain 0-based sequence number; renumber it after the shift.
        for ($r = 2; $r -le $lastRow; $r++) {
            $ws.Cells.Item($r, 1).Value = $r - 1
        }

        foreach ($rowNum in $fUpdates.Keys) {
            $ws.Cells.Item($rowNum, 6).Value = $fUpdates[$rowNum]
        }
    }
}
